$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: update Price (D) and Volume(1h) (E) columns.
# Numeric-looking price strings get a leading apostrophe (Excel's text-entry
# prefix) so they stay literal text -- matching the source's inline-string
# cells (e.g. thousands-grouped prices like '3.158.27') -- instead of being
# coerced to floating point numbers by Range.Value's auto-detection.

$ws.Range("D2").Value = '63.139.94'
$ws.Range("E2").Value = '  +1.10%  '
$ws.Range("D3").Value = '3.154.91'
$ws.Range("E3").Value = '  -0.92%  '
$ws.Range("D5").Value = '''597.28'
$ws.Range("E5").Value = '  +1.84%  '
$ws.Range("D6").Value = '''135.01'
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '3.154.35'
$ws.Range("E8").Value = '  -0.82%  '
$ws.Range("D9").Value = '''0.509'
$ws.Range("E9").Value = '  +1.47%  '
$ws.Range("E10").Value = '  +0.64%  '
$ws.Range("E11").Value = '  +2.35%  '
$ws.Range("D12").Value = '''0.451'
$ws.Range("E12").Value = '  +0.30%  '
$ws.Range("E13").Value = '  +1.61%  '
$ws.Range("D14").Value = '''34.63'
$ws.Range("E14").Value = '  +4.50%  '
$ws.Range("D15").Value = '3.675.99'
$ws.Range("E15").Value = '  -0.93%  '
$ws.Range("E16").Value = '  +1.62%  '
$ws.Range("D17").Value = '3.159.16'
$ws.Range("E17").Value = '  -1.20%  '
$ws.Range("D18").Value = '63.226.13'
$ws.Range("E18").Value = '  +1.27%  '
$ws.Range("D19").Value = '''6.54'
$ws.Range("E19").Value = '  -0.34%  '
$ws.Range("D20").Value = '''459.08'
$ws.Range("E20").Value = '  +0.69%  '
$ws.Range("D21").Value = '''13.89'
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").Value = '''0.693'
$ws.Range("E22").Value = '  -1.30%  '
$ws.Range("D23").Value = '''7.61'
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").Value = '''13.21'
$ws.Range("E24").Value = '  -1.41%  '
$ws.Range("D25").Value = '''82.86'
$ws.Range("E25").Value = '  +0.58%  '
$ws.Range("D26").Value = '''0.999'
$ws.Range("E26").Value = '  -0.14%  '
$ws.Range("D27").Value = '''2.68'
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("E29").Value = '  +2.87%  '
$ws.Range("D30").Value = '''7.66'
$ws.Range("E30").Value = '  -2.13%  '
$ws.Range("D31").Value = '''6.65'
$ws.Range("E31").Value = '  -3.83%  '
$ws.Range("D32").Value = '''26.97'
$ws.Range("E32").Value = '  -0.91%  '
$ws.Range("E33").Value = '  -1.55%  '
$ws.Range("D34").Value = '''2.39'
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").Value = '''1.01'
$ws.Range("E35").Value = '  -1.93%  '
$ws.Range("E36").Value = '  +1.62%  '
$ws.Range("D37").Value = '''51.14'
$ws.Range("E37").Value = '  -0.36%  '
$ws.Range("D38").Value = '0.0₃0721'
$ws.Range("E38").Value = '  +4.29%  '
$ws.Range("D39").Value = '''0.0388'
$ws.Range("E39").Value = '  +0.47%  '
$ws.Range("E40").Value = '  +1.42%  '
$ws.Range("E41").Value = '  -0.32%  '
$ws.Range("D42").Value = '''2.60'
$ws.Range("E42").Value = '  -0.89%  '
$ws.Range("D43").Value = '''389.18'
$ws.Range("E43").Value = '  -5.90%  '
$ws.Range("D44").Value = '2.775.83'
$ws.Range("E44").Value = '  -5.51%  '
$ws.Range("D45").Value = '''0.248'
$ws.Range("E45").Value = '  -0.53%  '
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").Value = '''35.68'
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("E48").Value = '  -2.33%  '
$ws.Range("D49").Value = '''125.72'
$ws.Range("E49").Value = '  +1.19%  '
$ws.Range("E50").Value = '  +0.59%  '
$ws.Range("D51").Value = '''24.83'
$ws.Range("E51").Value = '  -2.28%  '
